$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.74285042289024
$ws.Range("C2").Value = 10.68179576668495
$ws.Range("D2").Value = 12.1242983394225
$ws.Range("F2").Value = 26.75610915630249
$ws.Range("G2").Value = 22.73499991218451
$ws.Range("H2").Value = 12.82079565294901
$ws.Range("J2").Value = 11.26468044751778
$ws.Range("O2").Value = 18.70095410631783
$ws.Range("B3").Value = 14.94739067431882
$ws.Range("C3").Value = 10.05638762785308
$ws.Range("D3").Value = 12.01301686369503
$ws.Range("F3").Value = 26.83765529862294
$ws.Range("G3").Value = 22.88637762022889
$ws.Range("H3").Value = 12.89201821735609
$ws.Range("J3").Value = 11.2443736707799
$ws.Range("O3").Value = 18.82498018882316
$ws.Range("B4").Value = 14.43603411719885
$ws.Range("C4").Value = 9.65069987329476
$ws.Range("D4").Value = 11.94613701285682
$ws.Range("F4").Value = 26.89727908482472
$ws.Range("G4").Value = 22.99204908865799
$ws.Range("H4").Value = 12.938664330664
$ws.Range("J4").Value = 11.23419663554275
$ws.Range("O4").Value = 18.9072299129571
$ws.Range("B5").Value = 14.22207905388169
$ws.Range("C5").Value = 9.479998139963033
$ws.Range("D5").Value = 11.91927190494113
$ws.Range("F5").Value = 26.92396583301379
$ws.Range("G5").Value = 23.03827653945495
$ws.Range("H5").Value = 12.95840542180684
$ws.Range("J5").Value = 11.23062840469949
$ws.Range("O5").Value = 18.9422741660774
$ws.Range("B6").Value = 14.18622179170544
$ws.Range("C6").Value = 9.451330734637558
$ws.Range("D6").Value = 11.91483518408326
$ws.Range("F6").Value = 26.9285410764161
$ws.Range("G6").Value = 23.04614273993383
$ws.Range("H6").Value = 12.96172763396797
$ws.Range("J6").Value = 11.23007094504908
$ws.Range("O6").Value = 18.94818526930054
$ws.Range("B7").Value = 14.43317092537995
$ws.Range("C7").Value = 9.648419407711097
$ws.Range("D7").Value = 11.94577309316626
$ws.Range("F7").Value = 26.89762933412538
$ws.Range("G7").Value = 22.99265975736668
$ws.Range("H7").Value = 12.9389276011868
$ws.Range("J7").Value = 11.23414616546221
$ws.Range("O7").Value = 18.90769635835134
$ws.Range("B8").Value = 15.47344413553614
$ws.Range("C8").Value = 10.47069359767687
$ws.Range("D8").Value = 12.08564328204031
$ws.Range("F8").Value = 26.78223583914581
$ws.Range("G8").Value = 22.78453350810926
$ws.Range("H8").Value = 12.84474776884699
$ws.Range("J8").Value = 11.25720490734229
$ws.Range("O8").Value = 18.74244897503474
$ws.Range("B9").Value = 17.32443354234399
$ws.Range("C9").Value = 11.90876312452904
$ws.Range("D9").Value = 12.37028496797679
$ws.Range("F9").Value = 26.63228932244721
$ws.Range("G9").Value = 22.4788524519339
$ws.Range("H9").Value = 12.68322024509873
$ws.Range("J9").Value = 11.3204616036274
$ws.Range("O9").Value = 18.46706538221352
$ws.Range("B10").Value = 18.56175288437275
$ws.Range("C10").Value = 12.85675096504334
$ws.Range("D10").Value = 12.5841825767405
$ws.Range("F10").Value = 26.56932173742826
$ws.Range("G10").Value = 22.31868019904114
$ws.Range("H10").Value = 12.57870375826043
$ws.Range("J10").Value = 11.37771349960964
$ws.Range("O10").Value = 18.29480657540218
$ws.Range("B11").Value = 19.09690961782869
$ws.Range("C11").Value = 13.26416918252397
$ws.Range("D11").Value = 12.68219266394944
$ws.Range("F11").Value = 26.55104214738352
$ws.Range("G11").Value = 22.26018997653679
$ws.Range("H11").Value = 12.5342401821166
$ws.Range("J11").Value = 11.40604264779191
$ws.Range("O11").Value = 18.22305713151574
$ws.Range("B12").Value = 19.29550587021096
$ws.Range("C12").Value = 13.41500772676344
$ws.Range("D12").Value = 12.71938010281937
$ws.Range("F12").Value = 26.54561821723207
$ws.Range("G12").Value = 22.24013774225707
$ws.Range("H12").Value = 12.51784705748712
$ws.Range("J12").Value = 11.41709321493306
$ws.Range("O12").Value = 18.19684569794912
$ws.Range("B13").Value = 19.25291599660742
$ws.Range("C13").Value = 13.38267521649036
$ws.Range("D13").Value = 12.71136832160387
$ws.Range("F13").Value = 26.54671961796234
$ws.Range("G13").Value = 22.24436263983084
$ws.Range("H13").Value = 12.52135783242456
$ws.Range("J13").Value = 11.41469900950215
$ws.Range("O13").Value = 18.20244804459136
$ws.Range("B14").Value = 19.11332989939449
$ws.Range("C14").Value = 13.27664781075423
$ws.Range("D14").Value = 12.6852507937415
$ws.Range("F14").Value = 26.55056585261869
$ws.Range("G14").Value = 22.25849808049594
$ws.Range("H14").Value = 12.53288259934539
$ws.Range("J14").Value = 11.40694534640222
$ws.Range("O14").Value = 18.22088144629139
$ws.Range("B15").Value = 19.0272992901073
$ws.Range("C15").Value = 13.2112543588635
$ws.Range("D15").Value = 12.6692617358564
$ws.Range("F15").Value = 26.55311709292267
$ws.Range("G15").Value = 22.2674303605046
$ws.Range("H15").Value = 12.53999974353133
$ws.Range("J15").Value = 11.40223788442787
$ws.Range("O15").Value = 18.23229749599143
$ws.Range("B16").Value = 18.52621513162103
$ws.Range("C16").Value = 12.82964430468075
$ws.Range("D16").Value = 12.57778917252949
$ws.Range("F16").Value = 26.57072560599076
$ws.Range("G16").Value = 22.3227944935746
$ws.Range("H16").Value = 12.58167166235966
$ws.Range("J16").Value = 11.37590763032863
$ws.Range("O16").Value = 18.2996292789146
$ws.Range("B17").Value = 18.21166409188882
$ws.Range("C17").Value = 12.58942334722643
$ws.Range("D17").Value = 12.5218344657735
$ws.Range("F17").Value = 26.58418841516875
$ws.Range("G17").Value = 22.36046172381489
$ws.Range("H17").Value = 12.60802608300589
$ws.Range("J17").Value = 11.36033623413479
$ws.Range("O17").Value = 18.34263410323393
$ws.Range("B18").Value = 18.02813808963905
$ws.Range("C18").Value = 12.44901249138845
$ws.Range("D18").Value = 12.48971904621744
$ws.Range("F18").Value = 26.59290705949311
$ws.Range("G18").Value = 22.3834780866889
$ws.Range("H18").Value = 12.6234744217736
$ws.Range("J18").Value = 11.35159541759143
$ws.Range("O18").Value = 18.36799122116925
$ws.Range("B19").Value = 17.9655544386554
$ws.Range("C19").Value = 12.40108727894259
$ws.Range("D19").Value = 12.47885790967704
$ws.Range("F19").Value = 26.59602627603077
$ws.Range("G19").Value = 22.39150215428183
$ws.Range("H19").Value = 12.62875473278982
$ws.Range("J19").Value = 11.34867309381926
$ws.Range("O19").Value = 18.37668327988077
$ws.Range("B20").Value = 18.24541871704908
$ws.Range("C20").Value = 12.61522746086159
$ws.Range("D20").Value = 12.527784083474
$ws.Range("F20").Value = 26.58265428719207
$ws.Range("G20").Value = 22.35631191928691
$ws.Range("H20").Value = 12.60519059062273
$ws.Range("J20").Value = 11.36197157822958
$ws.Range("O20").Value = 18.33799174818335
$ws.Range("B21").Value = 19.15444029924176
$ws.Range("C21").Value = 13.30788413890882
$ws.Range("D21").Value = 12.69292038608062
$ws.Range("F21").Value = 26.54939540581021
$ws.Range("G21").Value = 22.25428902463667
$ws.Range("H21").Value = 12.52948542923581
$ws.Range("J21").Value = 11.40921406801346
$ws.Range("O21").Value = 18.21544103330851
$ws.Range("B22").Value = 19.7248734557419
$ws.Range("C22").Value = 13.74050984852675
$ws.Range("D22").Value = 12.80125977181293
$ws.Range("F22").Value = 26.53639382856357
$ws.Range("G22").Value = 22.19984511551512
$ws.Range("H22").Value = 12.48259782698388
$ws.Range("J22").Value = 11.44196914894236
$ws.Range("O22").Value = 18.14093770740219
$ws.Range("B23").Value = 19.42260800167869
$ws.Range("C23").Value = 13.51144912289672
$ws.Range("D23").Value = 12.74340856182318
$ws.Range("F23").Value = 26.54253155261746
$ws.Range("G23").Value = 22.22777405101241
$ws.Range("H23").Value = 12.50738523496234
$ws.Range("J23").Value = 11.42431715800032
$ws.Range("O23").Value = 18.18018741512357
$ws.Range("B24").Value = 18.23016661302989
$ws.Range("C24").Value = 12.60356859830469
$ws.Range("D24").Value = 12.52509409191301
$ws.Range("F24").Value = 26.58334481844379
$ws.Range("G24").Value = 22.3581838088391
$ws.Range("H24").Value = 12.60647159217717
$ws.Range("J24").Value = 11.36123158029799
$ws.Range("O24").Value = 18.3400885854651
$ws.Range("B25").Value = 16.84478058257225
$ws.Range("C25").Value = 11.5386716503947
$ws.Range("D25").Value = 12.29233100002786
$ws.Range("F25").Value = 26.66460864529614
$ws.Range("G25").Value = 22.55037236775467
$ws.Range("H25").Value = 12.72443450775921
$ws.Range("J25").Value = 11.30143863264845
$ws.Range("O25").Value = 18.53631322235229
